# Atualização de bases das ligas, do dia: 29-02-2024 às 07:50
#
# Several rows in "Mexico Liga MX Femenil" had their match-data columns
# (id .. PL_AhUnder, i.e. columns B and F through AC) mixed up between
# two (or three) adjacent rows. Column A (row index) and columns C/D/E
# (Div, Div Original Name, Date) are correct and must stay untouched;
# only the match-specific data needs to be swapped back into the right
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every column that carries match-specific data (everything except the
# sequential index in A and the Div/Div Original Name/Date in C:D:E).
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-RowData($ws, $cols, $r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

function Rotate-RowData($ws, $cols, $rA, $rB, $rC) {
    # rA takes rB's data, rB takes rC's data, rC takes rA's (original) data.
    foreach ($col in $cols) {
        $addrA = "$col$rA"
        $addrB = "$col$rB"
        $addrC = "$col$rC"
        $vA = $ws.Range($addrA).Value2
        $vB = $ws.Range($addrB).Value2
        $vC = $ws.Range($addrC).Value2
        $ws.Range($addrA).Value2 = $vB
        $ws.Range($addrB).Value2 = $vC
        $ws.Range($addrC).Value2 = $vA
    }
}

# Simple two-row swaps
Swap-RowData $ws $cols 28 29
Swap-RowData $ws $cols 47 48
Swap-RowData $ws $cols 101 102
Swap-RowData $ws $cols 109 110
Swap-RowData $ws $cols 131 132
Swap-RowData $ws $cols 133 134
Swap-RowData $ws $cols 149 150
Swap-RowData $ws $cols 221 222
Swap-RowData $ws $cols 232 233

# Three-row rotation: 229 <- 230 <- 231 <- 229
Rotate-RowData $ws $cols 229 230 231
